# Add two new columns "I0" (I) and "IF" (J) to the sheet, matching the
# header styling already used by the other header cells, and fill in
# the data rows with the value 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, borders, centered alignment)
# from the existing "IP" header cell (H1) onto the two new header cells
# so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2 and 3.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
